# "add dao insert and select" — add one new demo row to Customer, Drink and
# Food (the tables the new DAO insert/select methods exercise), converting
# the Drink/Food price & stock/sell columns from descriptive text to plain
# numbers, and leave the UI focused on the Food sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Customer: new row 6 — "wyw" (mirrors the existing "haichao" rows 4/5)
# ---------------------------------------------------------------------
$wsCustomer = $wb.Worksheets.Item("Customer")
$wsCustomer.Cells.Item(6, 1).Value = "'5"
$wsCustomer.Cells.Item(6, 2).Value = "wyw"
$wsCustomer.Cells.Item(6, 3).Value = "unknown"
$wsCustomer.Cells.Item(6, 4).Value = "'1"
$wsCustomer.Cells.Item(6, 5).Value = "'9"
$wsCustomer.Cells.Item(6, 6).Value = "'123456"
$wsCustomer.Cells.Item(6, 7).Value = "'1828888888"

# ---------------------------------------------------------------------
# Drink: prices / stock / sell amounts become plain numbers instead of
# "RMB nn" / "nn(库存)" / "nn(销量)" text, and a new drink is added.
# ---------------------------------------------------------------------
$wsDrink = $wb.Worksheets.Item("Drink")
$wsDrink.Cells.Item(2, 3).Value = 30
$wsDrink.Cells.Item(2, 4).Value = 100
$wsDrink.Cells.Item(2, 5).Value = 20

$wsDrink.Cells.Item(3, 3).Value = 32
$wsDrink.Cells.Item(3, 4).Value = 200
$wsDrink.Cells.Item(3, 5).Value = 40

$wsDrink.Cells.Item(4, 3).Value = 36
$wsDrink.Cells.Item(4, 4).Value = 180
$wsDrink.Cells.Item(4, 5).Value = 50

$wsDrink.Cells.Item(5, 1).Value = "'4"
$wsDrink.Cells.Item(5, 2).Value = "binghongcha"
$wsDrink.Cells.Item(5, 3).Value = "'30.0"
$wsDrink.Cells.Item(5, 4).Value = "'100"
$wsDrink.Cells.Item(5, 5).Value = "'1"

# ---------------------------------------------------------------------
# Food: same text->number conversion, plus a new "cookie" row.
# ---------------------------------------------------------------------
$wsFood = $wb.Worksheets.Item("Food")
$wsFood.Cells.Item(2, 3).Value = 10
$wsFood.Cells.Item(2, 4).Value = 200
$wsFood.Cells.Item(2, 5).Value = 30

$wsFood.Cells.Item(3, 3).Value = 20
$wsFood.Cells.Item(3, 4).Value = 150
$wsFood.Cells.Item(3, 5).Value = 50

$wsFood.Cells.Item(4, 1).Value = "'3"
$wsFood.Cells.Item(4, 2).Value = "cookie"
$wsFood.Cells.Item(4, 3).Value = "'30.0"
$wsFood.Cells.Item(4, 4).Value = "'100"
$wsFood.Cells.Item(4, 5).Value = "'1"

# ---------------------------------------------------------------------
# View state: restore each sheet's last selection, then leave the
# workbook focused on the Food sheet (it was "Customer" before).
# ---------------------------------------------------------------------
[void]$wsDrink.Range("E11").Select()
[void]$wsFood.Range("G17").Select()
[void]$wsFood.Activate()
